$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'36.130.68"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = "'  -4.36%  "
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = "'1.959.20"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = "'  -4.47%  "
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = "'  +0.21%  "
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = "'241.98"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = "'  -4.29%  "
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = "'0.620"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = "'  -4.99%  "
$ws.Range('E6').Style = "Normal"
$ws.Range('D7').Value = "'57.58"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = "'  -12.29%  "
$ws.Range('E7').Style = "Normal"
$ws.Range('E8').Value = "'  +0.18%  "
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').Value = "'0.369"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = "'  -2.58%  "
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = "'56.84"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = "'  -4.93%  "
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = "'0.0785"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = "'  +3.04%  "
$ws.Range('E11').Style = "Normal"
$ws.Range('E12').Value = "'  -0.96%  "
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = "'0.845"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = "'  -8.76%  "
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').Value = "'21.81"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = "'  +5.07%  "
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = "'13.88"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = "'  -8.92%  "
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = "'2.247.76"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = "'  -4.40%  "
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = "'5.36"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = "'  -3.69%  "
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = "'1.955.03"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = "'  -4.53%  "
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = "'36.068.96"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = "'  -4.23%  "
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = "'70.94"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = "'  -4.20%  "
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = "'0.0₃0848"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = "'  -3.63%  "
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = "'235.47"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = "'  -1.48%  "
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = "'5.17"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = "'  -3.44%  "
$ws.Range('E23').Style = "Normal"
$ws.Range('E24').Value = "'  -0.03%  "
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = "'2.52"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = "'  -5.97%  "
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = "'2.27"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = "'  -5.50%  "
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').Value = "'9.67"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = "'  +0.42%  "
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').Value = "'160.23"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = "'  -0.31%  "
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').Value = "'19.66"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = "'  -1.62%  "
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').Value = "'0.121"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = "'  +6.17%  "
$ws.Range('E30').Style = "Normal"
$ws.Range('E31').Value = "'  -2.61%  "
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = "'4.81"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = "'  -7.85%  "
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').Value = "'1.12"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = "'  -7.14%  "
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = "'0.0611"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = "'  -1.39%  "
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = "'4.37"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = "'  -8.00%  "
$ws.Range('E35').Style = "Normal"
$ws.Range('B36').Value = "'BinanceUSD"
$ws.Range('B36').Style = "Normal"
$ws.Range('C36').Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range('C36').Style = "Normal"
$ws.Range('D36').Value = "'1.00"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = "'  +0.24%  "
$ws.Range('E36').Style = "Normal"
$ws.Range('B37').Value = "'THORChain"
$ws.Range('B37').Style = "Normal"
$ws.Range('C37').Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range('C37').Style = "Normal"
$ws.Range('D37').Value = "'6.12"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = "'  +0.27%  "
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').Value = "'2.26"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = "'  -7.26%  "
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').Value = "'1.81"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = "'  -2.23%  "
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = "'3.03"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = "'  +6.14%  "
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').Value = "'0.0986"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = "'  -4.95%  "
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = "'1.21"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = "'  -2.01%  "
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').Value = "'2.88"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = "'  -1.62%  "
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = "'0.0212"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = "'  -3.55%  "
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').Value = "'1.08"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = "'  -5.17%  "
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').Value = "'91.21"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = "'  -4.32%  "
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = "'15.80"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = "'  -7.42%  "
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').Value = "'7.48"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = "'  -6.72%  "
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').Value = "'1.333.14"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = "'  -5.89%  "
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').Value = "'2.82"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = "'  -4.09%  "
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').Value = "'2.140.04"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = "'  -4.38%  "
$ws.Range('E51').Style = "Normal"
